$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "245.99"

# Row 4
Set-TextValue $ws.Range("D4") "5.361"

# Row 5
Set-TextValue $ws.Range("D5") "0.05734"

# Row 6
Set-TextValue $ws.Range("D6") "6.480"

# Row 7
Set-TextValue $ws.Range("D7") "3.144"

# Row 8
Set-TextValue $ws.Range("D8") "0.8186"

# Row 9
Set-TextValue $ws.Range("D9") "0.8730"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1380"
$ws.Range("E10").Value = "9WazirXWRX"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.06988"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

# Row 12
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D12") "0.03141"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.02940"
$ws.Range("E13").Value = "12BitrueCoinBTR"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09410"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# Row 15
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D15") "3.744"
$ws.Range("E15").Value = "14MCDexMCB"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D16") "0.001530"
$ws.Range("E16").Value = "15BitForexTokenBF"

# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D17") "0.04712"
$ws.Range("E17").Value = "16CoinExTokenCET"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D18") "0.0005980"
$ws.Range("E18").Value = "17OneONE"

# Row 19
Set-TextValue $ws.Range("D19") "0.006223"

# Row 20
Set-TextValue $ws.Range("D20") "0.001244"

# Row 21
Set-TextValue $ws.Range("D21") "0.004683"
$ws.Range("E21").Value = "20HotbitTokenHTBBestin24h"

# Row 22
Set-TextValue $ws.Range("D22") "0.00008802"

# Row 23
Set-TextValue $ws.Range("D23") "3.531"

# Row 24
Set-TextValue $ws.Range("D24") "2.141"

# Row 25
Set-TextValue $ws.Range("D25") "0.3183"

# Row 28
Set-TextValue $ws.Range("D28") "0.0002332"

# Row 40
Set-TextValue $ws.Range("D40") "0.03722"

# Row 41
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D41") "0.1060"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D42") "0.002601"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D43") "0.006461"
$ws.Range("E43").Value = "42KickTokenKICK"

# Row 44
Set-TextValue $ws.Range("D44") "0.008379"

# Row 45
Set-TextValue $ws.Range("D45") "0.00005259"

# Row 47
Set-TextValue $ws.Range("D47") "0.3999"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

# Row 48
Set-TextValue $ws.Range("D48") "0.002417"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

# Row 49
Set-TextValue $ws.Range("D49") "0.00002100"

# Row 50
Set-TextValue $ws.Range("D50") "0.0002000"
